# Add two new "Images aériennes" sub-items to the tag table on Sheet1,
# as children of the existing "images_aeriennes" tag (switch show tree
# on tag item page).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Extend the table (and therefore the sheet's used range / autofilter)
# by two rows.
$row42 = $tbl.ListRows.Add()
$row43 = $tbl.ListRows.Add()

# id / parent_id / name for the two new tag rows (filled column by
# column so new shared-string entries come out in the same order as
# the source workbook: id, id, parent_id(existing), parent_id(existing),
# name, name).
$ws.Range("A42").Value = "images_aeriennes_1"
$ws.Range("A43").Value = "images_aeriennes_2"

$ws.Range("B42").Value = "images_aeriennes"
$ws.Range("B43").Value = "images_aeriennes"

$ws.Range("C42").Value = "Images aériennes sous partie 1"
$ws.Range("C43").Value = "Images aériennes sous partie 2"

# doc_ids column stays empty but keeps the same wrapped-text style as
# the rest of the table.
$ws.Range("E42").WrapText = $true
$ws.Range("E43").WrapText = $true

# Row heights match the other single-line rows in the sheet.
$ws.Rows.Item(42).RowHeight = 16
$ws.Rows.Item(43).RowHeight = 16

# parent_id column is now a bit wider to fit "images_aeriennes".
$ws.Columns.Item(2).ColumnWidth = 13.5

# Restore the view to the top of the frozen pane / near the new rows,
# matching where the editor ended up after adding the rows.
$sel = $ws.Range("C44").Select()
$excel.ActiveWindow.ScrollRow = 2
